$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Topic: word embedding -> add Literature row
$ws.Rows(22).Insert()
$ws.Cells.Item(22, 1).Value = $ws.Cells.Item(21, 1).Value()
$ws.Cells.Item(22, 2).Value = $ws.Cells.Item(21, 2).Value()
$ws.Cells.Item(22, 3).Value = "Literature"
$ws.Cells.Item(22, 4).Value = $ws.Cells.Item(21, 4).Value()
$ws.Cells.Item(22, 5).Value = "As a literature researcher, I want to employ word embedding algorithms to identify and analyze shifts in literary themes and motifs across different periods and genres, uncovering evolutionary patterns in the use of language within the literary canon."

# Topic: voice recognition -> add Literature row
$ws.Rows(21).Insert()
$ws.Cells.Item(21, 1).Value = $ws.Cells.Item(20, 1).Value()
$ws.Cells.Item(21, 2).Value = $ws.Cells.Item(20, 2).Value()
$ws.Cells.Item(21, 3).Value = "Literature"
$ws.Cells.Item(21, 4).Value = $ws.Cells.Item(20, 4).Value()
$ws.Cells.Item(21, 5).Value = "As a literary critic, I want to use voice recognition technology to analyze audiobook recordings of classic novels and contemporary literature, extracting and comparing narrative styles and expressive techniques used by different authors."

# Topic: unsupervised clustering -> add Literature row
$ws.Rows(20).Insert()
$ws.Cells.Item(20, 1).Value = $ws.Cells.Item(19, 1).Value()
$ws.Cells.Item(20, 2).Value = $ws.Cells.Item(19, 2).Value()
$ws.Cells.Item(20, 3).Value = "Literature"
$ws.Cells.Item(20, 4).Value = $ws.Cells.Item(19, 4).Value()
$ws.Cells.Item(20, 5).Value = "As a literature researcher, I want to employ unsupervised clustering algorithms to analyze and categorize literary blogs and online discussions into thematic clusters such as literary criticism, author interviews, and reader responses, facilitating content analysis and trend identification in digital literary communities."

# Topic: text categorization -> add Literature row
$ws.Rows(19).Insert()
$ws.Cells.Item(19, 1).Value = $ws.Cells.Item(18, 1).Value()
$ws.Cells.Item(19, 2).Value = $ws.Cells.Item(18, 2).Value()
$ws.Cells.Item(19, 3).Value = "Literature"
$ws.Cells.Item(19, 4).Value = $ws.Cells.Item(18, 4).Value()
$ws.Cells.Item(19, 5).Value = "As a literary researcher, I want to use text categorization algorithms to classify novels into different genres automatically, based on textual features such as plot structure, character development, and thematic elements, to support genre-based analysis and recommendation systems."

# Topic: speech to text -> add Literature row
$ws.Rows(18).Insert()
$ws.Cells.Item(18, 1).Value = $ws.Cells.Item(17, 1).Value()
$ws.Cells.Item(18, 2).Value = $ws.Cells.Item(17, 2).Value()
$ws.Cells.Item(18, 3).Value = "Literature"
$ws.Cells.Item(18, 4).Value = $ws.Cells.Item(17, 4).Value()
$ws.Cells.Item(18, 5).Value = "As a literature educator, I want to integrate speech to text tools in educational settings to transcribe classroom discussions and oral presentations, facilitating documentation of student interactions and supporting inclusive learning environments for literary studies."

# Topic: sentiment analysis -> add Literature row
$ws.Rows(17).Insert()
$ws.Cells.Item(17, 1).Value = $ws.Cells.Item(16, 1).Value()
$ws.Cells.Item(17, 2).Value = $ws.Cells.Item(16, 2).Value()
$ws.Cells.Item(17, 3).Value = "Literature"
$ws.Cells.Item(17, 4).Value = $ws.Cells.Item(16, 4).Value()
$ws.Cells.Item(17, 5).Value = "As a literary critic, I want to perform sentiment analysis on book reviews to automatically classify the emotional tone (positive, negative, neutral) of reader feedback, providing insights into public reception and critical acclaim of literary works."

# Topic: semantic similarity -> add Literature row
$ws.Rows(16).Insert()
$ws.Cells.Item(16, 1).Value = $ws.Cells.Item(15, 1).Value()
$ws.Cells.Item(16, 2).Value = $ws.Cells.Item(15, 2).Value()
$ws.Cells.Item(16, 3).Value = "Literature"
$ws.Cells.Item(16, 4).Value = $ws.Cells.Item(15, 4).Value()
$ws.Cells.Item(16, 5).Value = "As a literary researcher, I want to develop semantic similarity metrics to quantify and compare the stylistic similarities between authors or literary movements, uncovering influences and stylistic evolution within the literary canon."

# Topic: random forest -> add Literature row
$ws.Rows(15).Insert()
$ws.Cells.Item(15, 1).Value = $ws.Cells.Item(14, 1).Value()
$ws.Cells.Item(15, 2).Value = $ws.Cells.Item(14, 2).Value()
$ws.Cells.Item(15, 3).Value = "Literature"
$ws.Cells.Item(15, 4).Value = $ws.Cells.Item(14, 4).Value()
$ws.Cells.Item(15, 5).Value = "As a literary researcher, I want to employ a random forest algorithm to classify literary texts into different literary movements or periods based on a combination of stylistic features such as vocabulary usage, sentence structure, and thematic elements."

# Topic: neural network -> add Literature row
$ws.Rows(14).Insert()
$ws.Cells.Item(14, 1).Value = $ws.Cells.Item(13, 1).Value()
$ws.Cells.Item(14, 2).Value = $ws.Cells.Item(13, 2).Value()
$ws.Cells.Item(14, 3).Value = "Literature"
$ws.Cells.Item(14, 4).Value = $ws.Cells.Item(13, 4).Value()
$ws.Cells.Item(14, 5).Value = "As a literature researcher, I want to develop a neural network architecture for authorship attribution tasks, accurately identifying the authorship of anonymous or disputed literary texts based on deep linguistic analysis and stylometric features."

# Topic: multi-label classification -> add Literature row
$ws.Rows(13).Insert()
$ws.Cells.Item(13, 1).Value = $ws.Cells.Item(12, 1).Value()
$ws.Cells.Item(13, 2).Value = $ws.Cells.Item(12, 2).Value()
$ws.Cells.Item(13, 3).Value = "Literature"
$ws.Cells.Item(13, 4).Value = $ws.Cells.Item(12, 4).Value()
$ws.Cells.Item(13, 5).Value = "As a literary researcher, I want to implement multi-label classification techniques to categorize novels into multiple genres simultaneously, recognizing hybrid genres or overlapping thematic elements in literary works."

# Topic: k-nearest neighbor -> add Literature row
$ws.Rows(12).Insert()
$ws.Cells.Item(12, 1).Value = $ws.Cells.Item(11, 1).Value()
$ws.Cells.Item(12, 2).Value = $ws.Cells.Item(11, 2).Value()
$ws.Cells.Item(12, 3).Value = "Literature"
$ws.Cells.Item(12, 4).Value = $ws.Cells.Item(11, 4).Value()
$ws.Cells.Item(12, 5).Value = "As a literature researcher, I want to employ k-NN clustering to analyze and categorize literary genres and subgenres based on shared characteristics such as narrative style, thematic content, and cultural influences, supporting systematic genre classification and literary taxonomy."

# Topic: keyword extraction -> add Literature row
$ws.Rows(11).Insert()
$ws.Cells.Item(11, 1).Value = $ws.Cells.Item(10, 1).Value()
$ws.Cells.Item(11, 2).Value = $ws.Cells.Item(10, 2).Value()
$ws.Cells.Item(11, 3).Value = "Literature"
$ws.Cells.Item(11, 4).Value = $ws.Cells.Item(10, 4).Value()
$ws.Cells.Item(11, 5).Value = "As a literary researcher, I want to use keyword extraction techniques to automatically identify and extract key themes and motifs from a collection of classic novels, facilitating thematic analysis and comparative literature studies."

# Topic: imbalanced dataset -> add Literature row
$ws.Rows(10).Insert()
$ws.Cells.Item(10, 1).Value = $ws.Cells.Item(9, 1).Value()
$ws.Cells.Item(10, 2).Value = $ws.Cells.Item(9, 2).Value()
$ws.Cells.Item(10, 3).Value = "Literature"
$ws.Cells.Item(10, 4).Value = $ws.Cells.Item(9, 4).Value()
$ws.Cells.Item(10, 5).Value = "As a literature curator, I want to manage imbalanced datasets in genre classification models, where niche or emerging literary genres are underrepresented compared to mainstream genres, to promote diversity and inclusivity in digital library collections and recommendations."

# Topic: feature selection -> add Literature row
$ws.Rows(9).Insert()
$ws.Cells.Item(9, 1).Value = $ws.Cells.Item(8, 1).Value()
$ws.Cells.Item(9, 2).Value = $ws.Cells.Item(8, 2).Value()
$ws.Cells.Item(9, 3).Value = "Literature"
$ws.Cells.Item(9, 4).Value = $ws.Cells.Item(8, 4).Value()
$ws.Cells.Item(9, 5).Value = "As a literary critic, I want to utilize feature selection techniques to identify and rank linguistic markers of literary movements or genres within a large corpus of novels or poems, supporting comprehensive genre analysis and literary classification."

# Topic: entity extraction -> add Literature row
$ws.Rows(8).Insert()
$ws.Cells.Item(8, 1).Value = $ws.Cells.Item(7, 1).Value()
$ws.Cells.Item(8, 2).Value = $ws.Cells.Item(7, 2).Value()
$ws.Cells.Item(8, 3).Value = "Literature"
$ws.Cells.Item(8, 4).Value = $ws.Cells.Item(7, 4).Value()
$ws.Cells.Item(8, 5).Value = "As a literature enthusiast, I want to use entity extraction techniques to automatically identify and extract character names, relationships, and attributes from complex narrative structures such as multi-volume novels or interconnected story collections, enhancing reader comprehension and engagement."

# Topic: document classification -> add Literature row
$ws.Rows(7).Insert()
$ws.Cells.Item(7, 1).Value = $ws.Cells.Item(6, 1).Value()
$ws.Cells.Item(7, 2).Value = $ws.Cells.Item(6, 2).Value()
$ws.Cells.Item(7, 3).Value = "Literature"
$ws.Cells.Item(7, 4).Value = $ws.Cells.Item(6, 4).Value()
$ws.Cells.Item(7, 5).Value = "As a literary critic, I want to utilize document classification algorithms to analyze and classify critical reviews and analyses of literary works into categories such as thematic interpretations, stylistic critiques, and cultural reflections, aiding in comprehensive literary criticism."

# Topic: decision tree -> add Literature row
$ws.Rows(6).Insert()
$ws.Cells.Item(6, 1).Value = $ws.Cells.Item(5, 1).Value()
$ws.Cells.Item(6, 2).Value = $ws.Cells.Item(5, 2).Value()
$ws.Cells.Item(6, 3).Value = "Literature"
$ws.Cells.Item(6, 4).Value = $ws.Cells.Item(5, 4).Value()
$ws.Cells.Item(6, 5).Value = "As a literary analyst, I want to employ a decision tree model to classify and categorize literary texts into thematic genres or periods based on their stylistic features, facilitating comprehensive literary analysis and comparative studies."

# Topic: conversational agent -> add Literature row
$ws.Rows(5).Insert()
$ws.Cells.Item(5, 1).Value = $ws.Cells.Item(4, 1).Value()
$ws.Cells.Item(5, 2).Value = $ws.Cells.Item(4, 2).Value()
$ws.Cells.Item(5, 3).Value = "Literature"
$ws.Cells.Item(5, 4).Value = $ws.Cells.Item(4, 4).Value()
$ws.Cells.Item(5, 5).Value = "As a literature enthusiast, I want to interact with a conversational agent that can discuss and recommend books based on my preferences and past readings, providing personalized literary suggestions and insights."

# Topic: cnn -> add Literature row
$ws.Rows(4).Insert()
$ws.Cells.Item(4, 1).Value = $ws.Cells.Item(3, 1).Value()
$ws.Cells.Item(4, 2).Value = $ws.Cells.Item(3, 2).Value()
$ws.Cells.Item(4, 3).Value = "Literature"
$ws.Cells.Item(4, 4).Value = $ws.Cells.Item(3, 4).Value()
$ws.Cells.Item(4, 5).Value = "As a literary scholar, I want to implement CNNs for text classification tasks to automatically categorize literary texts into genres or subgenres based on stylistic features and thematic content, facilitating comprehensive literary analysis."

# Topic: adversarial learning -> add Literature row
$ws.Rows(3).Insert()
$ws.Cells.Item(3, 1).Value = $ws.Cells.Item(2, 1).Value()
$ws.Cells.Item(3, 2).Value = $ws.Cells.Item(2, 2).Value()
$ws.Cells.Item(3, 3).Value = "Literature"
$ws.Cells.Item(3, 4).Value = $ws.Cells.Item(2, 4).Value()
$ws.Cells.Item(3, 5).Value = "As a literary critic, I want to apply adversarial learning techniques to analyze and detect stylistic forgeries or plagiarized passages in digital texts, ensuring the integrity and authenticity of literary works in online archives."

# Restore the active-cell selection to match the post-edit state
[void]$ws.Range("E8").Select()